# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G holds the "K" values (header in G1 is "K"). Update the
# per-row K values for rows 2-17 to reflect the regenerated save_data.
$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 1
    11 = 2
    13 = 0
    14 = 1
    15 = 2
    16 = 0
    17 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
